$wb = $excel.ActiveWorkbook

# Sheet "pokemon" gets new columns D:F with headers atk/def/pvMax and values
$wsPokemon = $wb.Worksheets.Item("pokemon")

$wsPokemon.Range("D1").Value = "atk"
$wsPokemon.Range("E1").Value = "def"
$wsPokemon.Range("F1").Value = "pvMax"

$wsPokemon.Range("D2").Value = 2
$wsPokemon.Range("E2").Value = 0
$wsPokemon.Range("F2").Value = 20

$wsPokemon.Range("D3").Value = 20
$wsPokemon.Range("E3").Value = 100
$wsPokemon.Range("F3").Value = 200

# Select a cell on the pokemon sheet so it becomes active / selected range matches diff
$wsPokemon.Range("G7").Select()

# Make pokemon sheet the active sheet (tab selected)
$wsPokemon.Activate()
